# Update Daily Report: 2026-02-02
# Appends the new day's (2026-01-30, Excel serial 46052) rows to Daily_Data,
# then refreshes the dependent Today_Summary and Monthly_Stats sheets to
# match the new totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Daily_Data: append rows 420..441 for date serial 46052
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily_Data")

# Columns: A=Date  B=Region_Type  C=PREV_TOTAL  D=RECEIVED  E=WITHDRAWN
#          F=NET_CHANGE  G=ADJUSTMENT  H=TOTAL_TODAY
$newRows = @(
    @(420, "ASAHI DEPOSITORY LLC Registered",                            0,          0, 0,          0,           0,  0),
    @(421, "ASAHI DEPOSITORY LLC Eligible",                              0,          0, 0,          0,           0,  0),
    @(422, "BRINK'S, INC. Registered",                            87949.747,        0, 0,          0, -11451.905,  76497.842),
    @(423, "BRINK'S, INC. Eligible",                              30578.352,        0, 0,          0,  11451.905,  42030.257),
    @(424, "CNT DEPOSITORY, INC. Registered",                       1246.06,        0, 0,          0,           0,   1246.06),
    @(425, "CNT DEPOSITORY, INC. Eligible",                               0,        0, 0,          0,           0,  0),
    @(426, "DELAWARE DEPOSITORY Registered",                       1633.941,        0, 0,          0,           0,   1633.941),
    @(427, "DELAWARE DEPOSITORY Eligible",                        18459.584,        0, 0,          0,           0,  18459.584),
    @(428, "HSBC BANK, USA Registered",                             1394.758,       0, 0,          0,           0,   1394.758),
    @(429, "HSBC BANK, USA Eligible",                        9281.978999999999,     0, 0,          0,           0, 9281.978999999999),
    @(430, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448,    0, 0,          0,           0,   2395.448),
    @(431, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible",          0,    0, 0,          0,           0,  0),
    @(432, "JP MORGAN CHASE BANK NA Registered",                  114985.579,       0, 0,          0,           0, 114985.579),
    @(433, "JP MORGAN CHASE BANK NA Eligible",                    125407.673,       0, 49923.162, -49923.162,    0,  75484.511),
    @(434, "LOOMIS INTERNATIONAL (US) LLC Registered",             63745.991,       0, 0,          0,           0,  63745.991),
    @(435, "LOOMIS INTERNATIONAL (US) LLC Eligible",               132077.206,      0, 0,          0,           0, 132077.206),
    @(436, "MALCA-AMIT USA, LLC Registered",                          395.145,      0, 0,          0,           0,    395.145),
    @(437, "MALCA-AMIT USA, LLC Eligible",                                  0,      0, 0,          0,           0,  0),
    @(438, "MANFRA, TORDELLA & BROOKES, LLC Registered",             50220.42,      0, 0,          0,           0,   50220.42),
    @(439, "MANFRA, TORDELLA & BROOKES, LLC Eligible",                1271.373,     0, 0,          0,           0,   1271.373),
    @(440, "STONEX PRECIOUS METALS LLC Registered",                  14122.765,     0, 0,          0,           0,  14122.765),
    @(441, "STONEX PRECIOUS METALS LLC Eligible",                        16.075,    0, 0,          0,           0,     16.075)
)

$dateSerial = 46052

foreach ($row in $newRows) {
    $r = $row[0]
    $daily.Cells.Item($r, 1).Value = $dateSerial
    $daily.Cells.Item($r, 2).Value = $row[1]
    $daily.Cells.Item($r, 3).Value = $row[2]
    $daily.Cells.Item($r, 4).Value = $row[3]
    $daily.Cells.Item($r, 5).Value = $row[4]
    $daily.Cells.Item($r, 6).Value = $row[5]
    $daily.Cells.Item($r, 7).Value = $row[6]
    $daily.Cells.Item($r, 8).Value = $row[7]
}

# Copy the date-column number format (style index used for column A) from the
# last pre-existing row so the new date cells are formatted consistently.
$daily.Range("A420:A441").NumberFormat = $daily.Range("A419").NumberFormat

# Dimension (A1:H441) is recalculated automatically by the host on save.

# ---------------------------------------------------------------------------
# 2) Today_Summary: refresh rows whose depository totals moved
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Today_Summary")

# Row 3: BRINK'S, INC.  -> Eligible (B), Registered (C)
$summary.Range("B3").Value = 42030.257
$summary.Range("C3").Value = 76497.842

# Row 8: JP MORGAN CHASE BANK NA -> Eligible (B), Total_Stock (D)
$summary.Range("B8").Value = 75484.511
$summary.Range("D8").Value = 190470.09

# ---------------------------------------------------------------------------
# 3) Monthly_Stats: refresh the 2026-01 month roll-up and the two
#    depository/type breakdown rows impacted by the new adjustments
# ---------------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")

# Row 2: 2026-01 month totals -> Eligible (B), Registered (C), Grand_Total (D)
$monthly.Range("B2").Value = 278620.985
$monthly.Range("C2").Value = 326637.949
$monthly.Range("D2").Value = 605258.934

# Row 9: BRINK'S, INC. Eligible -> TOTAL_TODAY (E)
$monthly.Range("E9").Value = 42030.257

# Row 10: BRINK'S, INC. Registered -> TOTAL_TODAY (E)
$monthly.Range("E10").Value = 76497.842

# Row 19: JP MORGAN CHASE BANK NA Eligible -> WITHDRAWN (D), TOTAL_TODAY (E)
$monthly.Range("D19").Value = 59929.312
$monthly.Range("E19").Value = 75484.511
